$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AF2").Value = "maa://25251 (93.02), ***maa://21730 (16.92), ***maa://39501 (17.65), *maa://36675 (60.0)"
$ws.Range("D6").Value = "maa://42407 (92.86)"
$ws.Range("P6").Value = "maa://31836 (89.47), maa://30381 (92.31)"
$ws.Range("X7").Value = "maa://22399 (94.89), *maa://22758 (71.93)"
$ws.Range("A8").Value = "更新日期：2024.11.10 13:16:41"
$ws.Range("AF9").Value = "maa://26206 (89.69), **maa://22865 (48.98)"
$ws.Range("D10").Value = "***maa://25695 (19.32), **maa://32237 (42.5), ***maa://34206 (18.18), ***maa://39951 (17.65), ***maa://39243 (28.57)"
$ws.Range("T10").Value = "maa://27395 (95.62), maa://22755 (87.27), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("X10").Value = "maa://22301 (97.47), maa://22726 (100.0)"
$ws.Range("D11").Value = "maa://36707 (99.35)"
$ws.Range("T11").Value = "maa://22747 (93.84), maa://22501 (98.18)"
$ws.Range("X11").Value = "maa://36713 (98.03)"
$ws.Range("AB12").Value = "maa://23669 (95.19), maa://36677 (93.33), maa://39872 (89.47)"
$ws.Range("AF12").Value = "*maa://28932 (77.78), *maa://20106 (63.96), *maa://22769 (64.29)"
$ws.Range("D13").Value = "maa://24999 (91.53), maa://36673 (92.42), maa://25001 (85.51)"
$ws.Range("H13").Value = "*maa://21248 (74.3), **maa://22728 (47.73)"
$ws.Range("P15").Value = "maa://24762 (89.93), *maa://22727 (70.0)"
$ws.Range("AF15").Value = "maa://21364 (80.74), *maa://22766 (70.75), *maa://36666 (77.33)"
$ws.Range("D16").Value = "maa://21441 (96.19), maa://36679 (92.31), maa://37650 (96.67)"
$ws.Range("D18").Value = "maa://24570 (96.88)"
$ws.Range("X18").Value = "maa://21917 (97.62), maa://22741 (83.33)"
$ws.Range("T19").Value = "maa://24386 (98.9)"
$ws.Range("L20").Value = "maa://41331 (84.06)"
$ws.Range("X21").Value = "maa://20110 (86.76), maa://34946 (91.89)"
$ws.Range("AB21").Value = "*maa://21443 (79.29), ***maa://23820 (29.82)"
$ws.Range("L23").Value = "maa://39756 (93.19), maa://39875 (94.74)"
$ws.Range("X24").Value = "maa://29988 (86.11), maa://23504 (92.88), **maa://22892 (39.86), *maa://25141 (77.42), maa://36663 (80.65), ***maa://22815 (23.08)"
$ws.Range("AF25").Value = "maa://20108 (96.21), maa://24621 (96.52), maa://36676 (96.3), maa://22771 (85.71), maa://37772 (100.0)"
$ws.Range("AB26").Value = "maa://42235 (90.57)"
$ws.Range("AF27").Value = "maa://24023 (96.92)"
$ws.Range("D28").Value = "maa://24465 (90.59), maa://25725 (83.13)"
$ws.Range("X28").Value = "maa://39929 (88.64), ***maa://39723 (14.29), maa://41749 (82.14)"
$ws.Range("AF28").Value = "maa://36660 (93.01), *maa://36701 (62.96)"
$ws.Range("AF29").Value = "*maa://24080 (69.17), ***maa://34960 (8.7), maa://42865 (90.0)"
$ws.Range("AB30").Value = "maa://42979 (95.74)"
$ws.Range("L31").Value = "maa://35926 (93.68), *maa://36258 (80.0)"
$ws.Range("T31").Value = "maa://30711 (96.49), maa://30768 (100.0)"
$ws.Range("T32").Value = "maa://41108 (89.36), maa://42859 (92.5), maa://41238 (95.45)"
$ws.Range("T34").Value = "maa://24526 (93.33)"
$ws.Range("AF38").Value = "maa://36697 (84.77)"
$ws.Range("H39").Value = "maa://25199 (85.32), maa://36670 (87.5), maa://30434 (87.72), ***maa://25036 (16.0)"
$ws.Range("H45").Value = "maa://21229 (85.56), maa://30807 (95.16), *maa://22767 (57.89), ***maa://20796 (13.79), *maa://42459 (60.0)"
$ws.Range("H46").Value = "maa://35931 (92.09)"
$ws.Range("H47").Value = "maa://27410 (95.92), maa://29661 (97.76), maa://28038 (84.62)"
$ws.Range("P49").Value = "*maa://39643 (68.42)"
$ws.Range("H55").Value = "maa://32532 (92.05)"
$ws.Range("H59").Value = "maa://27746 (83.5), maa://31270 (95.37)"
